$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the header text (PasteSpecial only copied formats, but make sure values stay correct)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-24
$data = @(
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(8, 9),
    @(5, 6),
    @(6, 7),
    @(6, 8),
    @(6, 6),
    @(4, 6),
    @(4, 5),
    @(4, 6),
    @(8, 8),
    @(8, 8),
    @(5, 7),
    @(6, 7),
    @(3, 5),
    @(6, 8),
    @(8, 9),
    @(6, 7),
    @(4, 5)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 9).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $r = $r + 1
}
